# "figured out mass allocation for crops"
# Adds a new "raw_name_edit" column (E) to the "2020" and "fish-dominant"
# sheets that consolidates/relabels some of the raw_name (D) categories.

$wb = $excel.ActiveWorkbook

# raw_name (column D) -> raw_name_edit (column E) lookup used on both sheets
$map = @{
    "soy protein concentrate" = "soy protein concentrate"
    "wheat gluten"            = "wheat gluten"
    "guar protein"            = "guar meal"
    "sunflower"               = "sunflower meal"
    "pea protein"             = "pea protein concentrate"
    "corn gluten"             = "corn gluten meal"
    "rapeseed oil"            = "canola oil"
    "linseed oil"             = "linseed oil"
    "soybean oil"             = "soy oil"
    "camelina oil"            = "canola oil"
    "coconut oil"             = "coconut oil"
    "wheat"                   = "wheat"
    "faba beans"              = "faba beans"
    "pea flour"               = "pea starch"
    "fish meal, forage fish"  = "fish meal, forage fish"
    "fish meal, cut offs"     = "fish meal, cut offs"
    "fish oil, forage fish"   = "fish oil, forage fish"
    "fish oil, cut offs"      = "fish oil, cut offs"
    "micro ingredients"       = "micro ingredients"
    "other"                   = "other"
    "soybean meal "           = "soybean meal "
}

foreach ($sheetName in @("2020", "fish-dominant")) {
    $ws = $wb.Worksheets.Item($sheetName)

    # header
    $ws.Range("E1").Value = "raw_name_edit"

    # last used row on this sheet (21 for 2020, 22 for fish-dominant)
    $lastRow = $ws.Cells.Item($ws.Rows.Count, 4).End(-4162).Row

    for ($r = 2; $r -le $lastRow; $r++) {
        $rawName = $ws.Cells.Item($r, 4).Text
        $ws.Cells.Item($r, 5).Value = $map[$rawName]
    }

    # match the column D formatting applied to the new column E
    $ws.Columns.Item(5).ColumnWidth = 13.33
}

# selections left behind by the editing session
$wb.Worksheets.Item("2020").Range("E16:E21").Select() | Out-Null
$wb.Worksheets.Item("fish-dominant").Range("M14").Select() | Out-Null
$wb.Worksheets.Item("fish-dominant").Activate() | Out-Null
